# Apply updated cryptocurrency market data to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.540.92"
$ws.Range("E2").Value = "  -0.24%  "

$ws.Range("D3").Value = "'2.289.62"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").Value = "'113.38"
$ws.Range("E5").Value = "  +16.99%  "

$ws.Range("D6").Value = "'268.27"
$ws.Range("E6").Value = "  +0.59%  "

$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("D10").Value = "'47.22"
$ws.Range("E10").Value = "  +3.07%  "

$ws.Range("D11").Value = "'0.0933"
$ws.Range("E11").Value = "  -0.34%  "

$ws.Range("D12").Value = "'8.47"
$ws.Range("E12").Value = "  +8.44%  "

$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").Value = "'15.54"
$ws.Range("E14").Value = "  +2.44%  "

$ws.Range("D15").Value = "'2.630.12"
$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").Value = "'0.843"
$ws.Range("E16").Value = "  -0.88%  "

$ws.Range("D17").Value = "'2.292.02"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").Value = "'43.535.76"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("E19").Value = "  +1.57%  "

$ws.Range("D20").Value = "'6.56"
$ws.Range("E20").Value = "  +6.37%  "

$ws.Range("D21").Value = "'72.27"
$ws.Range("E21").Value = "  +0.64%  "

$ws.Range("E22").Value = "  +3.09%  "

$ws.Range("D23").Value = "'232.45"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("E24").Value = "  +2.16%  "

$ws.Range("E25").Value = "  +13.46%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").Value = "'11.39"
$ws.Range("E27").Value = "  +1.97%  "

$ws.Range("D28").Value = "'42.06"
$ws.Range("E28").Value = "  +5.10%  "

$ws.Range("E29").Value = "  -1.70%  "

$ws.Range("E30").Value = "  +2.16%  "

$ws.Range("D31").Value = "'176.33"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("D32").Value = "'21.58"
$ws.Range("E32").Value = "  -3.29%  "

$ws.Range("D33").Value = "'0.0920"
$ws.Range("E33").Value = "  +4.16%  "

$ws.Range("D34").Value = "'5.49"
$ws.Range("E34").Value = "  +2.16%  "

$ws.Range("E35").Value = "  +0.40%  "

$ws.Range("D36").Value = "'4.68"
$ws.Range("E36").Value = "  +7.18%  "

$ws.Range("E37").Value = "  +0.40%  "

$ws.Range("E38").Value = "  -0.87%  "

$ws.Range("E39").Value = "  +10.09%  "

$ws.Range("D40").Value = "'2.42"
$ws.Range("E40").Value = "  +4.18%  "

$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "'73.43"
$ws.Range("E41").Value = "  +12.74%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.240"
$ws.Range("E42").Value = "  +0.75%  "

$ws.Range("D43").Value = "'13.59"
$ws.Range("E43").Value = "  +10.42%  "

$ws.Range("E44").Value = "  +5.01%  "

$ws.Range("E45").Value = "  +0.23%  "

$ws.Range("D46").Value = "'5.89"
$ws.Range("E46").Value = "  +12.55%  "

$ws.Range("D47").Value = "'8.72"
$ws.Range("E47").Value = "  -1.03%  "

$ws.Range("D48").Value = "'103.55"
$ws.Range("E48").Value = "  +5.46%  "

$ws.Range("D49").Value = "'0.100"
$ws.Range("E49").Value = "  -1.67%  "

$ws.Range("E50").Value = "  +2.81%  "

$ws.Range("E51").Value = "  +2.72%  "
